# Rename header row: "_old" -> "_FV2410", "_new" -> "_FV2504"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldHeaders = @("Segmentname_old","Segmentgruppe_old","Segment_old","Datenelement_old","Segment ID_old","Code_old","Qualifier_old","Beschreibung_old","Bedingungsausdruck_old","Bedingung_old")
$newHeaders = @("Segmentname_new","Segmentgruppe_new","Segment_new","Datenelement_new","Segment ID_new","Code_new","Qualifier_new","Beschreibung_new","Bedingungsausdruck_new","Bedingung_new")

for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $oldHeaders[$i] -replace "_old$", "_FV2410"
}

for ($i = 0; $i -lt $newHeaders.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = $newHeaders[$i] -replace "_new$", "_FV2504"
}

# Add frozen pane (freeze top row)
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Create table over the full used range
$range = $ws.Range("A1:U72")
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"
$table.TableStyle = ""
